# The dashboard's row 10 contains three "blog" widget cells whose text
# encodes a blog post series number ("ser: N"). This commit rotates each
# of those references forward by one post:
#   H10: ser 161 -> ser 162
#   D10: ser 162 -> ser 163
#   B10: ser 163 -> ser 164
# (i.e. the oldest referenced post, 161, is dropped and the next post,
# 164, takes its place at the end of the rotation.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 162"
$ws.Range("D10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 163"
$ws.Range("B10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 164"
